$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the mobile number values in column A (rows 2-5)
$ws.Range("A2").Value = 8015332963
$ws.Range("A3").Value = 8015332963
$ws.Range("A4").Value = 8015332963
$ws.Range("A5").Value = 8015332963

# Update the selection to B8
$ws.Range("B8").Select()
